$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("BonusPower", 2, 660, 136, 100, 200, "win"),
    @("BonusPower", 2, 1128, 127, 100, 200, "win"),
    @("BonusPower", 2.099999999999998, 828, 137, 100, 209.9999999999998, "win"),
    @("BonusPower", 2, 980, 146, 40, 80, "win")
)

$startRow = 24
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
